# Merge the "Random " + "Forest Regression " runs on the RESULT slide
# (slide 12, "Subtitle 2" shape, first paragraph) into a single run reading
# "Random Forest Regression ", matching a manual in-place text edit where a
# user simply typed to join two adjacent runs that share identical formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$target = "Random Forest Regression "
$full = $tr.Text
$idx = $full.IndexOf($target)

if ($idx -ge 0) {
    $start = $idx + 1
    $len = $target.Length
    $sub = $tr.Characters($start, $len)
    $sub.Text = $target
}
